# AP100 TestData workbook - remove the hard-coded environment URL,
# username and password that were stored in the "Input_Value" sheet
# (cells T2:V2) together with the hyperlink that pointed at that URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

$null = $ws.Activate()

# Clear the sensitive values (URL, UserName, Password) that used to live
# in T2, U2 and V2.
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()

# Remove the hyperlink that was attached to T2 (it pointed at the URL we
# just cleared).
$null = $ws.Range("T2").Hyperlinks.Delete()

# Reflect the selection the author left the sheet in.
$null = $ws.Range("T2:V2").Select()
